# Updates cryptocurrency price/volume figures in the cryptos list sheet.
# Mirrors the upstream GitHub Actions data refresh (prices + 1h volume %).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.028.27"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").Value = "2.337.96"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.31%  "

$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.510"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.79%  "

$ws.Range("E12").Value = "  -5.05%  "

$ws.Range("E13").Value = "  -1.40%  "

$ws.Range("D14").Value = "2.696.86"

$ws.Range("D15").Value = "2.342.16"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.802"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.65%  "

$ws.Range("D18").Value = "46.026.22"
$ws.Range("E18").Value = "  -1.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.45%  "

$ws.Range("D20").Value = "0.0₃0961"
$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("E21").Value = "  -3.55%  "

$ws.Range("E22").Value = "  -1.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.65%  "

$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.54%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "40.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.12%  "

$ws.Range("E28").Value = "  -3.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.35%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +14.13%  "

$ws.Range("E32").Value = "  +5.51%  "

$ws.Range("E33").Value = "  -6.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "144.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0771"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.111"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("E37").Value = "  -3.17%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.57%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.95%  "

$ws.Range("E41").Value = "  -3.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.94%  "

$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").Value = "1.840.06"
$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.69"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("E46").Value = "  -8.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.185"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "69.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.67%  "

$ws.Range("D49").Value = "2.566.86"
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "96.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.64%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.56%  "
